# Updates the cryptocurrency Price (D) and Volume(1h) (E) columns
# with freshly scraped figures, preserving each cell as literal text
# (matching the existing inline-string / text cells in the sheet).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rng = $ws.Range("D2:E2")
$rng.NumberFormat = "@"
$ws.Range("D2").Value = "307.97"
$ws.Range("E2").Value = "-2.46%"
$rng.ClearFormats()

$rng = $ws.Range("D3:E3")
$rng.NumberFormat = "@"
$ws.Range("D3").Value = "37.71"
$ws.Range("E3").Value = "-4.45%"
$rng.ClearFormats()

$rng = $ws.Range("E4")
$rng.NumberFormat = "@"
$ws.Range("E4").Value = "-0.22%"
$rng.ClearFormats()

$rng = $ws.Range("D5:E5")
$rng.NumberFormat = "@"
$ws.Range("D5").Value = "0.07873"
$ws.Range("E5").Value = "-3.77%"
$rng.ClearFormats()

$rng = $ws.Range("D6:E6")
$rng.NumberFormat = "@"
$ws.Range("D6").Value = "1.969"
$ws.Range("E6").Value = "-3.39%"
$rng.ClearFormats()

$rng = $ws.Range("D7:E7")
$rng.NumberFormat = "@"
$ws.Range("D7").Value = "4.346"
$ws.Range("E7").Value = "1.62%"
$rng.ClearFormats()

$rng = $ws.Range("D8:E8")
$rng.NumberFormat = "@"
$ws.Range("D8").Value = "8.248"
$ws.Range("E8").Value = "-0.10%"
$rng.ClearFormats()

$rng = $ws.Range("D9:E9")
$rng.NumberFormat = "@"
$ws.Range("D9").Value = "3.096"
$ws.Range("E9").Value = "-5.89%"
$rng.ClearFormats()

$rng = $ws.Range("D10:E10")
$rng.NumberFormat = "@"
$ws.Range("D10").Value = "0.9313"
$ws.Range("E10").Value = "-0.18%"
$rng.ClearFormats()

$rng = $ws.Range("D11:E11")
$rng.NumberFormat = "@"
$ws.Range("D11").Value = "0.1304"
$ws.Range("E11").Value = "-7.47%"
$rng.ClearFormats()

$rng = $ws.Range("D12:E12")
$rng.NumberFormat = "@"
$ws.Range("D12").Value = "0.1930"
$ws.Range("E12").Value = "-3.13%"
$rng.ClearFormats()

$rng = $ws.Range("D13:E13")
$rng.NumberFormat = "@"
$ws.Range("D13").Value = "0.08801"
$ws.Range("E13").Value = "-3.46%"
$rng.ClearFormats()

$rng = $ws.Range("D14:E14")
$rng.NumberFormat = "@"
$ws.Range("D14").Value = "0.03427"
$ws.Range("E14").Value = "-2.89%"
$rng.ClearFormats()

$rng = $ws.Range("D15:E15")
$rng.NumberFormat = "@"
$ws.Range("D15").Value = "0.09746"
$ws.Range("E15").Value = "-0.70%"
$rng.ClearFormats()

$rng = $ws.Range("D16:E16")
$rng.NumberFormat = "@"
$ws.Range("D16").Value = "0.001389"
$ws.Range("E16").Value = "-0.49%"
$rng.ClearFormats()

$rng = $ws.Range("D17:E17")
$rng.NumberFormat = "@"
$ws.Range("D17").Value = "0.005937"
$ws.Range("E17").Value = "-5.73%"
$rng.ClearFormats()

$rng = $ws.Range("E18")
$rng.NumberFormat = "@"
$ws.Range("E18").Value = "1,775.37%"
$rng.ClearFormats()

$rng = $ws.Range("E19")
$rng.NumberFormat = "@"
$ws.Range("E19").Value = "-2.56%"
$rng.ClearFormats()

$rng = $ws.Range("D20:E20")
$rng.NumberFormat = "@"
$ws.Range("D20").Value = "0.3434"
$ws.Range("E20").Value = "-0.74%"
$rng.ClearFormats()

$rng = $ws.Range("E21")
$rng.NumberFormat = "@"
$ws.Range("E21").Value = "-1.20%"
$rng.ClearFormats()

$rng = $ws.Range("D22:E22")
$rng.NumberFormat = "@"
$ws.Range("D22").Value = "4.996"
$ws.Range("E22").Value = "2.13%"
$rng.ClearFormats()

$rng = $ws.Range("D23:E23")
$rng.NumberFormat = "@"
$ws.Range("D23").Value = "0.2489"
$ws.Range("E23").Value = "1.48%"
$rng.ClearFormats()

$rng = $ws.Range("E24")
$rng.NumberFormat = "@"
$ws.Range("E24").Value = "-0.48%"
$rng.ClearFormats()

$rng = $ws.Range("D25:E25")
$rng.NumberFormat = "@"
$ws.Range("D25").Value = "0.001219"
$ws.Range("E25").Value = "-0.57%"
$rng.ClearFormats()

$rng = $ws.Range("D26:E26")
$rng.NumberFormat = "@"
$ws.Range("D26").Value = "0.004610"
$ws.Range("E26").Value = "-3.55%"
$rng.ClearFormats()

$rng = $ws.Range("E27")
$rng.NumberFormat = "@"
$ws.Range("E27").Value = "176.21%"
$rng.ClearFormats()

$rng = $ws.Range("D39:E39")
$rng.NumberFormat = "@"
$ws.Range("D39").Value = "0.02307"
$ws.Range("E39").Value = "3.49%"
$rng.ClearFormats()

$rng = $ws.Range("D40:E40")
$rng.NumberFormat = "@"
$ws.Range("D40").Value = "0.05033"
$ws.Range("E40").Value = "-4.25%"
$rng.ClearFormats()

$rng = $ws.Range("D41:E41")
$rng.NumberFormat = "@"
$ws.Range("D41").Value = "0.007516"
$ws.Range("E41").Value = "-0.40%"
$rng.ClearFormats()

$rng = $ws.Range("D42:E42")
$rng.NumberFormat = "@"
$ws.Range("D42").Value = "0.009754"
$ws.Range("E42").Value = "-0.12%"
$rng.ClearFormats()

$rng = $ws.Range("D43:E43")
$rng.NumberFormat = "@"
$ws.Range("D43").Value = "0.1357"
$ws.Range("E43").Value = "-1.35%"
$rng.ClearFormats()

$rng = $ws.Range("D44:E44")
$rng.NumberFormat = "@"
$ws.Range("D44").Value = "0.002088"
$ws.Range("E44").Value = "-2.95%"
$rng.ClearFormats()

$rng = $ws.Range("D45:E45")
$rng.NumberFormat = "@"
$ws.Range("D45").Value = "0.008004"
$ws.Range("E45").Value = "-15.72%"
$rng.ClearFormats()

$rng = $ws.Range("D46:E46")
$rng.NumberFormat = "@"
$ws.Range("D46").Value = "0.00006556"
$ws.Range("E46").Value = "1.56%"
$rng.ClearFormats()

$rng = $ws.Range("D47:E47")
$rng.NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000749"
$ws.Range("E47").Value = "-0.15%"
$rng.ClearFormats()

$rng = $ws.Range("D48:E48")
$rng.NumberFormat = "@"
$ws.Range("D48").Value = "0.002995"
$ws.Range("E48").Value = "8.17%"
$rng.ClearFormats()

$rng = $ws.Range("D50:E50")
$rng.NumberFormat = "@"
$ws.Range("D50").Value = "0.00002098"
$ws.Range("E50").Value = "-0.15%"
$rng.ClearFormats()

$rng = $ws.Range("D51:E51")
$rng.NumberFormat = "@"
$ws.Range("D51").Value = "0.0001998"
$ws.Range("E51").Value = "-0.15%"
$rng.ClearFormats()
